# Fixed Diagnosis, FileAssociation, FileFormat, FileType, NeuteredStatus, PrimeDiseaseSite
#
# The "CasesTab" row's Neo4j/Cypher query (column B) incorrectly returned an
# extra `Cohort` column. Remove the trailing `coalesce(co.cohort_description,
# '') AS `Cohort`` projection (and its now-dangling leading comma) from that
# query so it only returns the documented case-level fields.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newCasesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
MATCH (c)<--(diag:diagnosis)
MATCH (samp:sample)-->(c) 
  MATCH (f:file)-[*]->(c)
   WHERE f.file_type IN ["DNA Methylation Analysis File"] 
OPTIONAL MATCH (co:cohort)<-[*]-(c)
  WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`
'@

# Locate the row whose TabName (column A) is "CasesTab" rather than assuming
# a fixed row number, then overwrite its query cell (column B).
$casesTabCell = $ws.Range("A1:A1000").Find("CasesTab")
if ($casesTabCell -ne $null) {
    $row = $casesTabCell.Row
    $ws.Cells.Item($row, 2).Value = $newCasesQuery
} else {
    # Fallback to the known layout if the label can't be located.
    $ws.Range("B2").Value = $newCasesQuery
}
